$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 575; existing rows 575..648 shift down to 576..649.
$ws.Rows.Item(575).Insert()

# Populate the newly inserted row 575 with the new data record.
$ws.Cells.Item(575, 1).Value = 9
$ws.Cells.Item(575, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(575, 3).Value = "Metropolitana"
$ws.Cells.Item(575, 4).Value = 45142
$ws.Cells.Item(575, 5).Value = 13
$ws.Cells.Item(575, 6).Value = 100112012
$ws.Cells.Item(575, 7).Value = "Espinaca"
$ws.Cells.Item(575, 8).Value = "Sin especificar"
$ws.Cells.Item(575, 9).Value = "Primera"
$ws.Cells.Item(575, 10).Value = 160
$ws.Cells.Item(575, 11).Value = 6000
$ws.Cells.Item(575, 12).Value = 8000
$ws.Cells.Item(575, 13).Value = 7000
$ws.Cells.Item(575, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(575, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(575, 16).Value = 700
$ws.Cells.Item(575, 17).Value = 10
$ws.Cells.Item(575, 18).Value = "Hortaliza"
